$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.576.73"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.556.78"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'210.59"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'24.52"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D11").Value = "'0.0893"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "1.779.97"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "1.560.97"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "28.620.22"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "'61.29"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'229.75"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "0.0₃0671"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'8.95"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "'151.34"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'14.74"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "1.392.68"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "'1.04"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D41").Value = "'0.516"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "'64.00"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "1.692.67"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'0.869"
$ws.Range("E48").Value = "  -5.95%  "
$ws.Range("D49").Value = "'43.40"
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("D50").Value = "'85.08"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -0.66%  "
